$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.255.55'
$ws.Range("E2").Value = '  +0.54%  '

$ws.Range("D3").Value = '2.310.46'
$ws.Range("E3").Value = '  +0.57%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'301.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '

$ws.Range("D6").Value = "'98.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.04%  '

$ws.Range("D7").Value = "'0.520"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.24%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +0.95%  '

$ws.Range("D10").Value = "'36.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.36%  '

$ws.Range("D11").Value = "'0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("E12").Value = '  +0.42%  '

$ws.Range("D13").Value = "'17.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.92%  '

$ws.Range("D14").Value = "'6.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.99%  '

$ws.Range("D15").Value = '2.665.45'
$ws.Range("E15").Value = '  +0.07%  '

$ws.Range("D16").Value = '2.340.57'
$ws.Range("E16").Value = '  +2.84%  '

$ws.Range("D17").Value = "'0.793"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("D18").Value = '43.075.74'
$ws.Range("E18").Value = '  +0.36%  '

$ws.Range("D19").Value = "'13.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.42%  '

$ws.Range("E20").Value = '  +0.71%  '

$ws.Range("D21").Value = "'6.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").Value = "'68.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.70%  '

$ws.Range("D23").Value = "'238.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.99%  '

$ws.Range("D24").Value = "'2.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.67%  '

$ws.Range("D25").Value = "'0.990"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.33%  '

$ws.Range("E26").Value = '  -0.93%  '

$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").Value = "'25.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.79%  '

$ws.Range("D29").Value = "'166.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").Value = "'9.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("E31").Value = '  -6.61%  '

$ws.Range("D32").Value = "'33.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.51%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = "'5.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.35%  '

$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("D35").Value = "'18.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.40%  '

$ws.Range("D36").Value = "'4.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.76%  '

$ws.Range("E37").Value = '  -0.25%  '

$ws.Range("D38").Value = "'0.0693"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.20%  '

$ws.Range("E39").Value = '  +0.97%  '

$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("E41").Value = '  +1.40%  '

$ws.Range("D42").Value = "'2.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.69%  '

$ws.Range("D43").Value = '2.011.26'
$ws.Range("E43").Value = '  +1.34%  '

$ws.Range("E44").Value = '  -0.98%  '

$ws.Range("D45").Value = "'2.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.82%  '

$ws.Range("D46").Value = "'10.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.63%  '

$ws.Range("D47").Value = "'17.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.15%  '

$ws.Range("D48").Value = "'2.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.03%  '

$ws.Range("D49").Value = "'54.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.60%  '

$ws.Range("D50").Value = '2.536.83'
$ws.Range("E50").Value = '  +0.35%  '

$ws.Range("D51").Value = "'1.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.09%  '
